# Auto-generated edit script: updates column F ("想去人数" / want-to-go counts)
# across all 4 worksheets to match the target snapshot.
$wb = $excel.ActiveWorkbook

# --- Worksheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 988
$ws.Range("F3").Value = 571
$ws.Range("F4").Value = 8806
$ws.Range("F5").Value = 177
$ws.Range("F7").Value = 1949
$ws.Range("F8").Value = 6192
$ws.Range("F12").Value = 9123
$ws.Range("F13").Value = 10530
$ws.Range("F14").Value = 1199
$ws.Range("F15").Value = 1065
$ws.Range("F16").Value = 4806
$ws.Range("F17").Value = 755
$ws.Range("F18").Value = 403
$ws.Range("F20").Value = 315
$ws.Range("F21").Value = 172
$ws.Range("F22").Value = 1292
$ws.Range("F23").Value = 208
$ws.Range("F24").Value = 1852
$ws.Range("F25").Value = 841
$ws.Range("F26").Value = 1155
$ws.Range("F28").Value = 1986
$ws.Range("F29").Value = 394
$ws.Range("F30").Value = 573
$ws.Range("F31").Value = 2549
$ws.Range("F33").Value = 161
$ws.Range("F34").Value = 1640
$ws.Range("F35").Value = 87
$ws.Range("F37").Value = 401
$ws.Range("F38").Value = 885
$ws.Range("F39").Value = 564
$ws.Range("F40").Value = 3216
$ws.Range("F41").Value = 4201
$ws.Range("F42").Value = 227
$ws.Range("F44").Value = 480
$ws.Range("F45").Value = 557
$ws.Range("F47").Value = 886
$ws.Range("F48").Value = 222
$ws.Range("F49").Value = 4171

# --- Worksheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 8
$ws.Range("F6").Value = 10
$ws.Range("F8").Value = 30
$ws.Range("F22").Value = 65
$ws.Range("F26").Value = 57

# --- Worksheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 5692

# --- Worksheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 988
$ws.Range("F3").Value = 571
$ws.Range("F4").Value = 8806
$ws.Range("F5").Value = 177
$ws.Range("F7").Value = 8
$ws.Range("F8").Value = 6192
$ws.Range("F10").Value = 9123
$ws.Range("F11").Value = 9123
$ws.Range("F12").Value = 10530
$ws.Range("F13").Value = 30
$ws.Range("F14").Value = 1199
$ws.Range("F15").Value = 1065
$ws.Range("F16").Value = 4806
$ws.Range("F17").Value = 755
$ws.Range("F18").Value = 403
$ws.Range("F20").Value = 315
$ws.Range("F21").Value = 172
$ws.Range("F22").Value = 1292
$ws.Range("F23").Value = 208
$ws.Range("F24").Value = 841
$ws.Range("F25").Value = 1155
$ws.Range("F28").Value = 1986
$ws.Range("F29").Value = 394
$ws.Range("F30").Value = 573
$ws.Range("F31").Value = 2549
$ws.Range("F33").Value = 161
$ws.Range("F34").Value = 87
$ws.Range("F39").Value = 885
$ws.Range("F41").Value = 65
$ws.Range("F42").Value = 564
$ws.Range("F44").Value = 227
$ws.Range("F45").Value = 480
$ws.Range("F46").Value = 557
$ws.Range("F47").Value = 886
$ws.Range("F48").Value = 222
$ws.Range("F49").Value = 57

